$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Sel3 row (row 4): F4 10 -> 12, G4 8 -> 10
$ws.Range("F4").Value = 12
$ws.Range("G4").Value = 10

# Update active selection as shown in the diff (M13 -> J7)
$ws.Range("J7").Select()
